# Apply the "break out stock.yaml completed" edit:
#  - add a new "backup" column R
#  - backfill R2:R259 = 0 (two rows keep a non-zero backup marker: 183, 208)
#  - a handful of pre-existing "detect_structure" (Q) flags get cleared to 0
#  - row 259's isPivot (O) flag flips to 2
#  - 13 new weekly rows (260-272) are appended below the old last row (259)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New column header: R1 = "backup" (same bold/border/centered style as
#    the other header cells, so copy Q1's formatting across first).
# ---------------------------------------------------------------------------
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$ws.Range("R1").Value = "backup"

# ---------------------------------------------------------------------------
# 2. Backfill the new column for every existing data row (2-259) with 0,
#    then patch the two rows that actually carry a backup flag.
# ---------------------------------------------------------------------------
$ws.Range("R2:R259").Value = 0
$ws.Range("R183").Value = 2
$ws.Range("R208").Value = 2

# ---------------------------------------------------------------------------
# 3. A few existing rows had their detect_structure (Q) flag reset to 0.
# ---------------------------------------------------------------------------
$ws.Range("Q30").Value = 0
$ws.Range("Q38").Value = 0
$ws.Range("Q46").Value = 0
$ws.Range("Q56").Value = 0
$ws.Range("Q62").Value = 0

# ---------------------------------------------------------------------------
# 4. Row 259's isPivot (O) flag flips from 0 to 2.
# ---------------------------------------------------------------------------
$ws.Range("O259").Value = 2

# ---------------------------------------------------------------------------
# 5. Append 13 new weekly rows (260-272). Clone row 259's per-cell
#    formatting (date format on column A, plain on the rest) down first so
#    the new rows match the existing look, then fill in the values.
# ---------------------------------------------------------------------------
$newRows = @(
    @(45453, 2569.5,             2624.85009765625,  2520,               2577.550048828125, 2577.550048828125, 592237,  2024, 6, 10, 0, 0, 0, 24, 0, 2, 2),
    @(45460, 2577.550048828125,  2717.25,            2535.75,            2655.10009765625,  2655.10009765625,  1115224, 2024, 6, 17, 0, 0, 0, 25, 0, 0, 0),
    @(45467, 2663.949951171875,  2734.199951171875,  2643.25,            2680.85009765625,  2680.85009765625,  926049,  2024, 6, 24, 0, 0, 0, 26, 0, 0, 0),
    @(45474, 2681,               2742.75,            2638.199951171875,  2698.800048828125, 2698.800048828125, 642604,  2024, 7, 1,  0, 0, 0, 27, 0, 0, 0),
    @(45481, 2712,               2846.39990234375,   2631.449951171875,  2812,               2812,               1349257, 2024, 7, 8,  0, 0, 0, 28, 0, 0, 0),
    @(45488, 2820.35009765625,   2953.949951171875,  2780.85009765625,   2798.199951171875, 2798.199951171875, 1757147, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 1),
    @(45495, 2789.949951171875,  3051,               2728.050048828125,  3032.949951171875, 3032.949951171875, 1575475, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(45502, 3064,               3198.39990234375,   2765,               2788.75,            2788.75,            3388995, 2024, 7, 29, 0, 0, 0, 31, 1, 0, 0),
    @(45509, 2716.64990234375,   2848.39990234375,   2626.25,            2777.25,            2777.25,            1580013, 2024, 8, 5,  0, 0, 0, 32, 0, 0, 0),
    @(45516, 2780,               2821.300048828125,  2662,               2748.699951171875, 2748.699951171875, 665195,  2024, 8, 12, 0, 0, 0, 33, 0, 0, 0),
    @(45523, 2770,               2992.89990234375,   2749,               2896.64990234375,  2896.64990234375,  1576807, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0),
    @(45530, 2898,               3067.550048828125,  2898,               3000.050048828125, 3000.050048828125, 1612292, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0),
    @(45537, 3009,               3122,               2939.800048828125,  2983.800048828125, 2983.800048828125, 1841692, 2024, 9, 2,  0, 0, 0, 36, 0, 0, 0)
)

$targetRow = 260
foreach ($rowValues in $newRows) {
    $ws.Range("A259:Q259").Copy()
    $ws.Range("A" + $targetRow + ":Q" + $targetRow).PasteSpecial(-4122)

    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($targetRow, $col).Value = $rowValues[$col - 1]
    }

    $targetRow++
}

$excel.CutCopyMode = 0
